$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: after the "Integrate with Networking" bullet, add two
# new bullets:
#   - "Create methods to handle game invites to friends" (ilvl 0)
#   - "Implement methods to handle invites accepted and rejected" (ilvl 1)
# -----------------------------------------------------------------

$targetIdx = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "Integrate with Networking*") {
        $targetIdx = $i
    }
}

$target = $d.Paragraphs.Item($targetIdx)
$target.Range.InsertParagraphAfter()

$p1 = $d.Paragraphs.Item($targetIdx + 1)
$p1.Range.Text = "Create methods to handle game invites to friends"
$p1.Range.ListFormat.ListLevelNumber = 1

$p1.Range.InsertParagraphAfter()

$p2 = $d.Paragraphs.Item($targetIdx + 2)
$p2.Range.Text = "Implement methods to handle invites accepted and rejected"
$p2.Range.ListFormat.ListLevelNumber = 2

# -----------------------------------------------------------------
# Change 2: merge the "Leaves a / 3 hour / buffer..." runs (which
# had a proofErr gramStart/gramEnd pair around "3 hour") into one
# continuous run of text.
# -----------------------------------------------------------------

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    "Leaves a 3 hour buffer in case something goes wrong. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Leaves a 3 hour buffer in case something goes wrong. ",
    2
) | Out-Null
